$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ESTIMATE")

# Row labels 18-27: re-worded charge/line-item descriptions
$ws.Range("A25").Value = "Gross Azure Cluster Charge"
$ws.Range("A26").Value = "Gross OpenShift Cluster Licensing Charge"
$ws.Range("A24").Value = "Cluster Azure Other 2 Charge"
$ws.Range("A23").Value = "Cluster Azure Other 1 Charge"
$ws.Range("A22").Value = "Cluster Azure Load Balancer Charge"
$ws.Range("A21").Value = "Cluster Azure Storage (~40% of Compute)"
$ws.Range("A20").Value = "Cluster Azure Master Node Compute"
$ws.Range("A19").Value = "Cluster Azure App Node Compute"
$ws.Range("A18").Value = "Cluster OpenShift App Node Licensing"
$ws.Range("A27").Value = "Gross Total Cluster Charge"

# Header: "Annualized" -> "Annual"
$ws.Range("N1").Value = "Annual"

# Row 21 (Cluster Azure Storage) used to be a hard-coded estimate; it is now
# derived as 40% of the combined App Node + Master Node compute charges.
$ws.Range("B21:M21").Formula = "=40%*SUM(B19:B20)"

# Azure Credit (row 28) - no monthly credit assumed by default any more.
$ws.Range("B28").Value = 0

# EA / volume discount (row 30) - reset the flat 15% assumption to 0.
$ws.Range("B30:M30").Value = 0

# Column A needed to widen to fit the longer relabeled charge descriptions.
$ws.Columns.Item(1).ColumnWidth = 48.6

# Restore the active selection Excel had when it was last saved.
$ws.Range("M30").Select()
